$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.111.88"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.892.16"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'307.45"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5144"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("D8").Value = "'0.3743"
$ws.Range("E8").Value = "  +3.32%  "
$ws.Range("D9").Value = "'0.07214"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").Value = "'21.21"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").Value = "'0.9060"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "'0.07643"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").Value = "1.889.79"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "'95.13"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").Value = "'5.277"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'0.000008486"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'14.47"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("D20").Value = "27.124.27"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'5.074"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").Value = "2.115.07"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("D24").Value = "'6.409"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'145.78"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").Value = "'1.787"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "'2.233"
$ws.Range("E27").Value = "  +8.46%  "
$ws.Range("D28").Value = "'18.11"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").Value = "'114.60"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "'4.965"
$ws.Range("E30").Value = "  +6.10%  "
$ws.Range("D31").Value = "'4.844"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").Value = "'0.09192"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").Value = "'0.05094"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Value = "'1.238"
$ws.Range("E34").Value = "  +7.71%  "
$ws.Range("D35").Value = "'0.7745"
$ws.Range("E35").Value = "  +4.09%  "
$ws.Range("D36").Value = "'2.991"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "'3.290"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").Value = "'2.634"
$ws.Range("E38").Value = "  +5.01%  "
$ws.Range("D39").Value = "'0.01999"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'0.5594"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "'6.668"
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("D43").Value = "'8.981"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("D44").Value = "'117.77"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'0.1512"
$ws.Range("E45").Value = "  +3.01%  "
$ws.Range("D46").Value = "'0.4803"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").Value = "'10.24"
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("D48").Value = "'1.0000"
$ws.Range("D49").Value = "'1.597"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").Value = "'37.64"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").Value = "'64.07"
$ws.Range("E51").Value = "  +1.71%  "
